$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddress, $textValue)
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $textValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "315.66"
Set-TextValue "E2" "3.31%"
Set-TextValue "E3" "-1.28%"
Set-TextValue "D4" "5.123"
Set-TextValue "E4" "0.42%"
Set-TextValue "D5" "0.08105"
Set-TextValue "E5" "2.83%"
Set-TextValue "D6" "2.135"
Set-TextValue "E6" "-0.32%"
Set-TextValue "D7" "8.001"
Set-TextValue "E7" "0.91%"
Set-TextValue "D8" "4.149"
Set-TextValue "E8" "1.02%"
Set-TextValue "D9" "0.9261"
Set-TextValue "E9" "0.32%"
Set-TextValue "D10" "0.1014"
Set-TextValue "E10" "4.28%"
Set-TextValue "D11" "0.1872"
Set-TextValue "E11" "0.79%"
Set-TextValue "D12" "0.09225"
Set-TextValue "E12" "6.35%"
Set-TextValue "D13" "0.03604"
Set-TextValue "E13" "1.26%"
Set-TextValue "D14" "0.09904"
Set-TextValue "E14" "-0.38%"
Set-TextValue "D15" "0.001436"
Set-TextValue "E15" "0.53%"
Set-TextValue "D16" "0.005671"
Set-TextValue "E16" "0.84%"
Set-TextValue "E17" "0.32%"
Set-TextValue "D18" "2.831"
Set-TextValue "E18" "7.42%"
Set-TextValue "D19" "0.3368"
Set-TextValue "E19" "-0.80%"
Set-TextValue "D20" "0.1330"
Set-TextValue "E20" "0.98%"
Set-TextValue "D21" "5.154"
Set-TextValue "E21" "-0.53%"
Set-TextValue "D22" "0.2220"
Set-TextValue "E22" "0.79%"
Set-TextValue "D23" "0.04569"
Set-TextValue "E23" "0.10%"
Set-TextValue "D24" "0.001247"
Set-TextValue "E24" "0.93%"
Set-TextValue "D25" "0.004707"
Set-TextValue "E25" "-6.84%"
Set-TextValue "D26" "0.0001252"
Set-TextValue "E26" "-21.91%"
Set-TextValue "D27" "0.0004505"
Set-TextValue "E27" "-5.05%"
Set-TextValue "D39" "0.01957"
Set-TextValue "E39" "5.99%"
Set-TextValue "D40" "0.04869"
Set-TextValue "E40" "2.03%"
Set-TextValue "D41" "0.007736"
Set-TextValue "E41" "3.14%"
Set-TextValue "E42" "-0.66%"
Set-TextValue "D43" "0.007836"
Set-TextValue "D44" "0.002143"
Set-TextValue "E44" "-3.78%"
Set-TextValue "D45" "0.01163"
Set-TextValue "E45" "5.42%"
Set-TextValue "D46" "0.00006535"
Set-TextValue "E46" "3.47%"
Set-TextValue "E47" "0.21%"
Set-TextValue "D48" "39.20"
Set-TextValue "E48" "-17.47%"
Set-TextValue "D49" "0.001702"
Set-TextValue "E49" "-14.82%"
Set-TextValue "E50" "0.21%"
Set-TextValue "E51" "0.21%"
